# Karel/Eclipse handout update ("Updated karel handout to reference
# eclipse online instructions"):
#  1. "...described in Handout #5" -> "...described on the course website"
#  2. Move the _GoBack bookmark from the "that linked / page" sentence to
#     just before the "This button is the " paragraph.
#  3. "...bring up the Eclipse window shown on the last page of Handout #5"
#     -> "...bring up the window shown at the end of the online Eclipse
#     instructions"
#  4. Remove the stray empty "Separator" paragraph right after the
#     "... section of the Eclipse screen:" paragraph.

$d = $word.ActiveDocument

# --- 1. First "Handout #5" reference ---------------------------------------
$found1 = $d.Content.Find.Execute(
    "Once you have downloaded a copy of Eclipse as described in Handout #5",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Once you have downloaded a copy of Eclipse as described on the course website",
    2)
if (-not $found1) {
    Write-Output "WARNING: first Handout #5 reference not found"
}

# --- 3. Second "Handout #5" reference ---------------------------------------
$found2 = $d.Content.Find.Execute(
    "From here, your next step is to start up Eclipse, which will bring up the Eclipse window shown on the last page of Handout #5",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "From here, your next step is to start up Eclipse, which will bring up the window shown at the end of the online Eclipse instructions",
    2)
if (-not $found2) {
    Write-Output "WARNING: second Handout #5 reference not found"
}

# --- 2. Move the _GoBack bookmark -------------------------------------------
# Find the paragraph that starts with "This button is the " and collapse a
# range to its very start, then (re)plant the _GoBack bookmark there. Adding
# a bookmark named "_GoBack" implicitly removes any existing one elsewhere in
# the document, which is exactly the "move" the diff describes.
$movedBookmark = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("This button is the ")) {
        $target = $para.Range.Duplicate()
        $target.Collapse(1)
        $d.Bookmarks.Add("_GoBack", $target) | Out-Null
        $movedBookmark = $true
        break
    }
}
if (-not $movedBookmark) {
    Write-Output "WARNING: 'This button is the ' paragraph not found; bookmark not moved"
}

# --- 4. Remove the extra empty Separator paragraph --------------------------
$removedSeparator = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.Contains("section of the Eclipse screen:")) {
        $next = $d.Paragraphs($i + 1)
        if ($next.Style.NameLocal -eq "Separator" -and $next.Range.Text.Trim() -eq "") {
            $next.Range.Delete() | Out-Null
            $removedSeparator = $true
        }
        break
    }
}
if (-not $removedSeparator) {
    Write-Output "WARNING: separator paragraph after 'section of the Eclipse screen:' not removed"
}
